# Update "想去人数" (F) and "最低票价" (G) figures per upstream data refresh
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 10333
$ws1.Range("F5").Value = 763
$ws1.Range("F8").Value = 482
$ws1.Range("F9").Value = 447
$ws1.Range("F11").Value = 276
$ws1.Range("F13").Value = 12928
$ws1.Range("G13").Value = 238
$ws1.Range("F14").Value = 12928
$ws1.Range("G14").Value = 238
$ws1.Range("F16").Value = 54
$ws1.Range("F18").Value = 198
$ws1.Range("F19").Value = 149
$ws1.Range("F20").Value = 190
$ws1.Range("F21").Value = 2773
$ws1.Range("F24").Value = 2116
$ws1.Range("F25").Value = 125
$ws1.Range("F27").Value = 417
$ws1.Range("F29").Value = 2173
$ws1.Range("F30").Value = 1144
$ws1.Range("F31").Value = 4324
$ws1.Range("F33").Value = 3900
$ws1.Range("F34").Value = 944
$ws1.Range("F35").Value = 2677
$ws1.Range("F36").Value = 3106
$ws1.Range("F37").Value = 96
$ws1.Range("F38").Value = 1401
$ws1.Range("F39").Value = 218
$ws1.Range("F40").Value = 798
$ws1.Range("F41").Value = 60
$ws1.Range("F42").Value = 161
$ws1.Range("F43").Value = 597
$ws1.Range("F44").Value = 849
$ws1.Range("F46").Value = 175
$ws1.Range("F47").Value = 350
$ws1.Range("F48").Value = 129
$ws1.Range("F49").Value = 194
$ws1.Range("F50").Value = 216

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F19").Value = 41

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 10333
$ws4.Range("F5").Value = 763
$ws4.Range("F7").Value = 482
$ws4.Range("F8").Value = 447
$ws4.Range("F10").Value = 276
$ws4.Range("F11").Value = 12928
$ws4.Range("G11").Value = 238
$ws4.Range("F12").Value = 12928
$ws4.Range("G12").Value = 238
$ws4.Range("F15").Value = 54
$ws4.Range("F17").Value = 198
$ws4.Range("F19").Value = 2773
$ws4.Range("F21").Value = 2116
$ws4.Range("F22").Value = 125
$ws4.Range("F24").Value = 417
$ws4.Range("F26").Value = 2173
$ws4.Range("F27").Value = 1144
$ws4.Range("F31").Value = 4324
$ws4.Range("F32").Value = 3900
$ws4.Range("F33").Value = 944
$ws4.Range("F34").Value = 2677
$ws4.Range("F35").Value = 3106
$ws4.Range("F36").Value = 96
$ws4.Range("F39").Value = 218
$ws4.Range("F40").Value = 798
$ws4.Range("F41").Value = 60
$ws4.Range("F42").Value = 597
$ws4.Range("F43").Value = 41
$ws4.Range("F44").Value = 849
$ws4.Range("F46").Value = 175
$ws4.Range("F47").Value = 350
$ws4.Range("F48").Value = 129
$ws4.Range("F49").Value = 194
$ws4.Range("F50").Value = 216
